# Append new Q-Learning / Policy-Iteration / Value-Iteration / Human-player
# match results to the results log sheet (rows 14-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2025-06-25 19:34:11", "Policy Iteration", "{'win': 1, 'loss': 1, 'draw': 0}"),
    @("2025-06-25 19:34:20", "Policy Iteration", "{'win': 1, 'loss': 2, 'draw': 1}"),
    @("2025-06-25 19:34:22", "Policy Iteration", "{'win': 3, 'loss': 2, 'draw': 1}"),
    @("2025-06-25 19:35:09", "Policy Iteration", "{'win': 0, 'loss': 0, 'draw': 2}"),
    @("2025-06-25 19:35:13", "Policy Iteration", "{'win': 0, 'loss': 1, 'draw': 3}"),
    @("2025-06-25 19:35:17", "Policy Iteration", "{'win': 0, 'loss': 1, 'draw': 5}"),
    @("2025-06-25 19:47:15", "Value Iteration",  "{'win': 0, 'loss': 1, 'draw': 1}"),
    @("2025-06-25 19:47:20", "Value Iteration",  "{'win': 1, 'loss': 2, 'draw': 1}"),
    @("2025-06-25 19:47:23", "Value Iteration",  "{'win': 3, 'loss': 2, 'draw': 1}"),
    @("2025-06-25 19:47:28", "Value Iteration",  "{'win': 3, 'loss': 2, 'draw': 3}"),
    @("2025-06-25 19:47:31", "Value Iteration",  "{'win': 4, 'loss': 2, 'draw': 4}"),
    @("2025-06-25 19:48:24", "Joueur Humain",    "{'win': 1, 'loss': 1, 'draw': 0}")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]

    # Stat/Hyper columns (D-I) stay empty for these rows, same as the
    # other rows logged without hyper-parameters (e.g. rows 2-4, 9-13).
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 9).Value = ""
}
